$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.55
$ws.Range("G2").Value = 4.7
$ws.Range("H2").Value = 2.08
$ws.Range("I2").Value = 2.4
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 3.5
$ws.Range("L2").Value = 1.38
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 2.76
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 1.61
$ws.Range("Q2").Value = 2.24
$ws.Range("R2").Value = 1.24
$ws.Range("S2").Value = 3.8
$ws.Range("T2").Value = 1.89
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.71
$ws.Range("W2").Value = 1.28
$ws.Range("F3").Value = 1.1
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 50
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 1.05
$ws.Range("AA4").Value = 470
$ws.Range("AE4").Value = 110
$ws.Range("AG4").Value = 12
$ws.Range("AK4").Value = 40
$ws.Range("AN4").Value = 13.5
$ws.Range("AO4").Value = 60
$ws.Range("H4").Value = 3.1
$ws.Range("AA5").Value = 1000
$ws.Range("AC5").Value = 7.8
$ws.Range("AD5").Value = 17
$ws.Range("AF5").Value = 14
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 21
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AO5").Value = 85
$ws.Range("G5").Value = 2.38
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 4.1
$ws.Range("L5").Value = 1.38
$ws.Range("N5").Value = 3.05
$ws.Range("O5").Value = 1.38
$ws.Range("Q5").Value = 2.12
$ws.Range("S5").Value = 3.9
$ws.Range("V5").Value = 1.33
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 12
$ws.Range("Z5").Value = 65
$ws.Range("AA6").Value = 870
$ws.Range("AB6").Value = 8.6
$ws.Range("AC6").Value = 14.5
$ws.Range("AD6").Value = 55
$ws.Range("AE6").Value = 310
$ws.Range("AF6").Value = 7.2
$ws.Range("AG6").Value = 11.5
$ws.Range("AH6").Value = 40
$ws.Range("AI6").Value = 240
$ws.Range("AJ6").Value = 9.4
$ws.Range("AK6").Value = 14.5
$ws.Range("AL6").Value = 46
$ws.Range("AN6").Value = 4.9
$ws.Range("AO6").Value = 450
$ws.Range("G6").Value = 1.29
$ws.Range("H6").Value = 14.5
$ws.Range("I6").Value = 16.5
$ws.Range("J6").Value = 6.4
$ws.Range("K6").Value = 6.8
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 5.1
$ws.Range("O6").Value = 1.23
$ws.Range("P6").Value = 2.34
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 1.52
$ws.Range("S6").Value = 2.78
$ws.Range("T6").Value = 2.3
$ws.Range("V6").Value = 1.06
$ws.Range("W6").Value = 4.5
$ws.Range("X6").Value = 22
$ws.Range("Y6").Value = 55
$ws.Range("Z6").Value = 150
$ws.Range("AB7").Value = 8.2
$ws.Range("AC7").Value = 8.2
$ws.Range("AD7").Value = 19
$ws.Range("AF7").Value = 13
$ws.Range("AG7").Value = 11.5
$ws.Range("AH7").Value = 22
$ws.Range("AI7").Value = 290
$ws.Range("AJ7").Value = 28
$ws.Range("AK7").Value = 65
$ws.Range("AN7").Value = 22
$ws.Range("F7").Value = 1.98
$ws.Range("G7").Value = 2.2
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4.8
$ws.Range("J7").Value = 3.15
$ws.Range("K7").Value = 3.7
$ws.Range("L7").Value = 1.47
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 2.92
$ws.Range("O7").Value = 1.42
$ws.Range("P7").Value = 1.66
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.24
$ws.Range("S7").Value = 4.2
$ws.Range("T7").Value = 1.93
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.26
$ws.Range("W7").Value = 1.83
$ws.Range("X7").Value = 11.5
$ws.Range("Y7").Value = 14
$ws.Range("Z7").Value = 34
$ws.Range("AA8").Value = 270
$ws.Range("AB8").Value = 8.8
$ws.Range("AC8").Value = 7.6
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 150
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 21
$ws.Range("AK8").Value = 80
$ws.Range("AN8").Value = 27
$ws.Range("F8").Value = 2.28
$ws.Range("G8").Value = 2.44
$ws.Range("H8").Value = 3.6
$ws.Range("K8").Value = 3.35
$ws.Range("L8").Value = 1.49
$ws.Range("N8").Value = 2.98
$ws.Range("P8").Value = 1.66
$ws.Range("Q8").Value = 2.28
$ws.Range("R8").Value = 1.23
$ws.Range("S8").Value = 4.4
$ws.Range("W8").Value = 1.69
$ws.Range("X8").Value = 11
$ws.Range("Z8").Value = 30
$ws.Range("AD9").Value = 110
$ws.Range("AF9").Value = 8.2
$ws.Range("AJ9").Value = 8.2
$ws.Range("AM9").Value = 480
$ws.Range("AN9").Value = 3.95
$ws.Range("I9").Value = 29
$ws.Range("M9").Value = 1.02
$ws.Range("P9").Value = 2.38
$ws.Range("Q9").Value = 1.59
$ws.Range("R9").Value = 1.56
$ws.Range("T9").Value = 2.62
$ws.Range("U9").Value = 1.5
$ws.Range("X9").Value = 32
$ws.Range("Z9").Value = 320
$ws.Range("AB10").Value = 8.4
$ws.Range("AC10").Value = 15.5
$ws.Range("AF10").Value = 7.6
$ws.Range("AJ10").Value = 9.6
$ws.Range("AK10").Value = 18
$ws.Range("AM10").Value = 340
$ws.Range("AN10").Value = 5.4
$ws.Range("F10").Value = 1.29
$ws.Range("G10").Value = 1.32
$ws.Range("H10").Value = 14
$ws.Range("L10").Value = 1.33
$ws.Range("M10").Value = 1.04
$ws.Range("N10").Value = 4.4
$ws.Range("P10").Value = 2.16
$ws.Range("Q10").Value = 1.69
$ws.Range("S10").Value = 2.78
$ws.Range("X10").Value = 22
$ws.Range("Z10").Value = 190
$ws.Range("AB11").Value = 8.2
$ws.Range("AC11").Value = 8.8
$ws.Range("AF11").Value = 11.5
$ws.Range("AG11").Value = 11
$ws.Range("AI11").Value = 100
$ws.Range("AL11").Value = 50
$ws.Range("AM11").Value = 180
$ws.Range("F11").Value = 1.88
$ws.Range("G11").Value = 1.94
$ws.Range("K11").Value = 3.75
$ws.Range("S11").Value = 3.85
$ws.Range("V11").Value = 1.25
$ws.Range("W11").Value = 2.06
$ws.Range("X11").Value = 15
$ws.Range("Y11").Value = 15.5
$ws.Range("Z11").Value = 980
$ws.Range("AB12").Value = 8
$ws.Range("AI12").Value = 90
$ws.Range("AJ12").Value = 24
$ws.Range("AN12").Value = 21
$ws.Range("AO12").Value = 100
$ws.Range("G12").Value = 2.06
$ws.Range("K12").Value = 3.6
$ws.Range("N12").Value = 2.98
$ws.Range("Q12").Value = 2.2
$ws.Range("T12").Value = 1.97
$ws.Range("U12").Value = 1.86
$ws.Range("W12").Value = 1.94
$ws.Range("X12").Value = 11
$ws.Range("AB13").Value = 14.5
$ws.Range("AD13").Value = 170
$ws.Range("AF13").Value = 9
$ws.Range("AG13").Value = 22
$ws.Range("AJ13").Value = 7.8
$ws.Range("AK13").Value = 18
$ws.Range("AM13").Value = 600
$ws.Range("AN13").Value = 2.32
$ws.Range("F13").Value = 1.07
$ws.Range("H13").Value = 44
$ws.Range("I13").Value = 55
$ws.Range("N13").Value = 8.6
$ws.Range("P13").Value = 3.6
$ws.Range("S13").Value = 1.84
$ws.Range("T13").Value = 2.74
$ws.Range("U13").Value = 1.46
$ws.Range("Y13").Value = 160
$ws.Range("H14").Value = 2.04
$ws.Range("I14").Value = 2.16
$ws.Range("P14").Value = 1.67
$ws.Range("Q14").Value = 2.04
